$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = "-"
$ws.Range("C2").Value = "-"
$ws.Range("C3").Value = "Desenho Técnico"
$ws.Range("C6").Value = "Desenho Técnico"
